# Restructure ontology: remove mfd_hab1=Urban if mfd_areatype=Urban
#
# For every data row whose mfd_hab1 (column N) is "Urban":
#   - habitat_typenumber (column F): 6410 -> 4100, 6420 -> 4200
#   - mfd_hab1 (column N): "Urban" -> "Greenspaces" (old mfd_hab2 value)
#   - mfd_hab2 (column O): takes the old mfd_hab3 (column P) value
#   - mfd_hab3 (column P): removed (cell cleared)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$habitatNumberMap = @{ "6410" = "4100"; "6420" = "4200" }

$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {

    $hab1 = $ws.Cells.Item($row, 14).Value()   # column N: mfd_hab1

    if ($hab1 -eq "Urban") {

        # --- column F: habitat_typenumber ---
        $rawNumber = $ws.Cells.Item($row, 6).Value()
        $numberKey = [string]$rawNumber
        if ($habitatNumberMap.ContainsKey($numberKey)) {
            $newNumber = $habitatNumberMap[$numberKey]
            $fCell = $ws.Cells.Item($row, 6)
            $fCell.NumberFormat = "@"
            $fCell.Value = $newNumber
            $fCell.Style = "Normal"
        }

        # --- columns N, O, P: shift hab2 -> hab1, hab3 -> hab2, clear hab3 ---
        $hab2 = $ws.Cells.Item($row, 15).Value()   # column O: mfd_hab2
        $hab3 = $ws.Cells.Item($row, 16).Value()   # column P: mfd_hab3

        $ws.Cells.Item($row, 14).Value = $hab2     # N = old O
        $ws.Cells.Item($row, 15).Value = $hab3     # O = old P
        $ws.Cells.Item($row, 16).ClearContents()   # P removed
    }
}
